$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '72.831.39'
$ws.Range('E2').Value = '  +2.50%  '
$ws.Range('D3').Value = '3.981.17'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''611.99'
$ws.Range('E5').Value = '  +14.11%  '
$ws.Range('D6').Value = '''163.55'
$ws.Range('E6').Value = '  +10.66%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '''0.755'
$ws.Range('E9').Value = '  +2.07%  '
$ws.Range('E10').Value = '  +1.64%  '
$ws.Range('D11').Value = '''54.40'
$ws.Range('E11').Value = '  -1.29%  '
$ws.Range('E12').Value = '  +1.01%  '
$ws.Range('E13').Value = '  +3.52%  '
$ws.Range('D14').Value = '4.617.44'
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('D15').Value = '3.992.86'
$ws.Range('E15').Value = '  +1.54%  '
$ws.Range('D16').Value = '''1.26'
$ws.Range('E16').Value = '  +8.51%  '
$ws.Range('D17').Value = '''14.15'
$ws.Range('E17').Value = '  +1.48%  '
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('E19').Value = '  +0.45%  '
$ws.Range('D20').Value = '72.563.38'
$ws.Range('E20').Value = '  +2.36%  '
$ws.Range('D21').Value = '''440.65'
$ws.Range('E21').Value = '  +4.54%  '
$ws.Range('D23').Value = '''96.67'
$ws.Range('E23').Value = '  -0.73%  '
$ws.Range('D24').Value = '''3.47'
$ws.Range('E24').Value = '  -3.35%  '
$ws.Range('D25').Value = '''14.43'
$ws.Range('E25').Value = '  -0.83%  '
$ws.Range('D26').Value = '''4.32'
$ws.Range('E26').Value = '  +11.96%  '
$ws.Range('D27').Value = '''11.49'
$ws.Range('E27').Value = '  +2.09%  '
$ws.Range('E28').Value = '  +1.10%  '
$ws.Range('D29').Value = '''10.55'
$ws.Range('E29').Value = '  -1.86%  '
$ws.Range('D30').Value = '''36.50'
$ws.Range('E30').Value = '  +0.19%  '
$ws.Range('D31').Value = '''7.80'
$ws.Range('E31').Value = '  +0.73%  '
$ws.Range('D32').Value = '''13.95'
$ws.Range('E32').Value = '  +4.79%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').Value = '''48.75'
$ws.Range('E34').Value = '  -4.80%  '
$ws.Range('D35').Value = '''71.46'
$ws.Range('E35').Value = '  +8.83%  '
$ws.Range('D36').Value = '''659.32'
$ws.Range('E36').Value = '  -2.56%  '
$ws.Range('D37').Value = '0.0₃0902'
$ws.Range('E37').Value = '  +11.14%  '
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('E39').Value = '  -1.10%  '
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('D41').Value = '''0.999'
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('D42').Value = '''3.35'
$ws.Range('E42').Value = '  +5.18%  '
$ws.Range('E43').Value = '  +0.38%  '
$ws.Range('D44').Value = '''0.0490'
$ws.Range('E44').Value = '  +1.77%  '
$ws.Range('D45').Value = '''10.63'
$ws.Range('E45').Value = '  +4.57%  '
$ws.Range('E46').Value = '  +0.95%  '
$ws.Range('D47').Value = '''2.66'
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('D48').Value = '''3.39'
$ws.Range('E48').Value = '  +1.66%  '
$ws.Range('D49').Value = '2.904.29'
$ws.Range('E49').Value = '  +12.33%  '
$ws.Range('E50').Value = '  +2.30%  '
$ws.Range('E51').Value = '  +4.80%  '
